# Update gh-pages generated output (commit 456a3b4):
#  - bump a handful of "interested" counters (column F) on rows 3,4,6,7
#  - insert a new expo entry (2024.04.21 合肥·银魂only) as the new row 8,
#    pushing the existing "合肥·梦时空SPO1动漫展" row down to row 9
# Applies to both the "展览" and "全部类型" sheets (they carry duplicate data).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- counter bumps -----------------------------------------------
    $ws.Range("F3").Value = 2447
    $ws.Range("F4").Value = 459
    $ws.Range("F6").Value = 6520
    $ws.Range("F7").Value = 354

    # --- insert a new row 8, pushing the old row 8 down to row 9 -----
    # Move the existing row 8 contents down to row 9 (values + formats)
    # without letting Excel fabricate a brand-new style record the way
    # Rows.Insert() does.
    $ws.Range("A8:I8").Copy()
    $ws.Range("A9").PasteSpecial(-4104)
    $ws.Application.CutCopyMode = $false

    # Row 9's A cell needs the same bold/bordered style as column A
    # elsewhere; copy it explicitly (PasteSpecial above left A9 using the
    # default style) and restore the running index value to 8.
    $ws.Range("A7").Copy()
    $ws.Range("A9").PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false
    $ws.Range("A9").Value = 8

    # New row 8: "合肥·银魂only" entry. A8 keeps its existing value/style
    # (7) — only B8..I8 are populated with the new expo's data.
    # B8's text ("2024.04.21") looks like a bare date, which Excel would
    # otherwise silently convert to a date serial; force text, assign,
    # then drop the now-unneeded number format back to General so the
    # cell ends up styled exactly like its plain-text neighbours.
    $ws.Range("B8").NumberFormat = "@"
    $ws.Range("B8").Value = "2024.04.21"
    $ws.Range("B8").ClearFormats()
    $ws.Range("C8").Value = "合肥·银魂only"
    $ws.Range("D8").Value = "濉溪路118号 合肥栢景假日酒店"
    $ws.Range("E8").Value = "2024.04.21 09:00-04.21 17:00"
    $ws.Range("F8").Value = 0
    $ws.Range("G8").Value = "不可售"
    $ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=82145"
    $ws.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202402/A0Tb5SQ51709091316985.jpeg"
}
